$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fully clear column B (contents + formatting) since it is removed entirely
$ws.Range("B1:B3").Clear()

# Clear contents (but keep formatting, e.g. the bold/centered header style) for column A
$ws.Range("A1:A3").ClearContents()

# New values for column A, rows 1 (header) through 28
$values = @(
    "leiden_fusion",
    "MeV.2.1",
    "MeV.2.8",
    "MeV.1.4.2",
    "MeV.4.21",
    "MeV.1.4.5",
    "MeV.1.4.7",
    "MeV.1.4.15",
    "MeV.1.4.6",
    "MeV.1.4.4",
    "MeV.1.4.20",
    "MeV.1.4.1",
    "MeV.1.4.11",
    "MeV.1.4.8",
    "MeV.4.12",
    "MeV.4.4",
    "MeV.1.4.0",
    "MeV.3.17",
    "MeV.4.31",
    "MeV.4.1",
    "MeV.4.34",
    "MeV.1.4.13",
    "MeV.3.30",
    "MeV.4.26",
    "MeV.1.4.12",
    "MeV.1.4.21",
    "MeV.4.30",
    "MeV.NA"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
